$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 2: new person (Juan Miguel David Becerra Tobar, cod_rh=0001458832)
$ws.Range("B2").Value = "Juan Miguel David"
$ws.Range("D2").Value = "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001458832"
$ws.Range("C2").Value = "Becerra Tobar"

# Update selection / view state
$ws.Range("C3").Select()
